# Apply edits to the "Dashboard" sheet (4th sheet) of the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# Set the new "Y" values on D2/E2 and D10/E10 (previously blank cells).
$ws.Range("D2").Value = "Y"
$ws.Range("E2").Value = "Y"
$ws.Range("D10").Value = "Y"
$ws.Range("E10").Value = "Y"

# Highlight D5/E5, D8/E8, D9/E9 with a yellow fill (style moves from s=2 to s=1).
$ws.Range("D5:E5").Interior.Color = 65535
$ws.Range("D8:E8").Interior.Color = 65535
$ws.Range("D9:E9").Interior.Color = 65535

# Update the active selection to E10, matching the new sheetView selection.
$ws.Range("E10").Select()
